$d = $word.ActiveDocument

# Each "old" answer string below is unique within the document, so a
# whole-content Find/Replace (wdReplaceAll = 2) safely retargets only the
# single matching table cell for every pair.
$replacements = @(
    @("921÷4=230, 1", "921÷7=131, 4"),
    @("701÷3=233, 2", "663÷4=165, 3"),
    @("975÷6=162, 3", "726÷8=90, 6"),
    @("316÷2=158, 0", "664÷6=110, 4"),
    @("234÷6=39, 0",  "404÷9=44, 8"),
    @("355÷7=50, 5",  "322÷6=53, 4"),
    @("143÷7=20, 3",  "131÷5=26, 1"),
    @("493÷2=246, 1", "497÷4=124, 1"),
    @("252÷9=28, 0",  "982÷2=491, 0"),
    @("912÷2=456, 0", "317÷3=105, 2"),
    @("970÷6=161, 4", "634÷9=70, 4"),
    @("684÷7=97, 5",  "849÷4=212, 1"),
    @("893÷4=223, 1", "206÷3=68, 2"),
    @("687÷6=114, 3", "579÷5=115, 4"),
    @("406÷7=58, 0",  "108÷6=18, 0"),
    @("515÷7=73, 4",  "825÷9=91, 6"),
    @("704÷2=352, 0", "980÷5=196, 0"),
    @("907÷3=302, 1", "340÷4=85, 0"),
    @("892÷7=127, 3", "105÷8=13, 1"),
    @("254÷5=50, 4",  "435÷4=108, 3"),
    @("900÷7=128, 4", "225÷3=75, 0"),
    @("278÷8=34, 6",  "795÷7=113, 4"),
    @("563÷9=62, 5",  "680÷9=75, 5"),
    @("978÷3=326, 0", "940÷8=117, 4"),
    @("286÷8=35, 6",  "146÷2=73, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: could not find text to replace: $old"
    }
}

Write-Output "Done: applied $($replacements.Count) replacements."
